$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.441.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.821.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5417"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4039"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07679"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.119"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.322"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.629"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.001"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.822.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001088"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06606"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.059"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.448.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.271"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.465"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.037.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1110"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.679"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07371"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.643"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2240"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02339"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.206"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6283"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.181"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.403"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.699"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5850"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.004"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.92%  "
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06868"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.75%  "
